# Update countries & provincias Spain
# Daily refresh of the COVID country-stats table on sheet "Pais".
# A handful of countries changed rank (sorted descending by "Casos
# totales" in column B), which shuffles which country name sits on a
# given row; the rest just get refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 19:56"

# --- Rows whose country name shifted because of re-sorting ----------
# Marruecos overtook Belgica & Ecuador
$ws.Range("A32").Value = "Marruecos"
$ws.Range("A33").Value = "Belgica"
$ws.Range("A34").Value = "Ecuador"

# Irlanda overtook Libia & Kenia
$ws.Range("A72").Value = "Irlanda"
$ws.Range("A73").Value = "Libia"
$ws.Range("A74").Value = "Kenia"

# --- Refreshed statistics (Casos totales, Nuevos casos, Casos
#     activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----

# Estados Unidos
$ws.Range("B4").Value = 7918732
$ws.Range("C4").Value = 24254
$ws.Range("D4").Value = 5071257
$ws.Range("E4").Value = 2628543
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 285
$ws.Range("H4").Value = 218932

# India
$ws.Range("B5").Value = 7037329
$ws.Range("C5").Value = 60321
$ws.Range("D5").Value = 6045564
$ws.Range("E5").Value = 883719
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 596
$ws.Range("H5").Value = 108046

# Brasil
$ws.Range("B6").Value = 5073483
$ws.Range("C6").Value = 16293
$ws.Range("D6").Value = 4433595
$ws.Range("E6").Value = 489865
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 331
$ws.Range("H6").Value = 150023

# Irak
$ws.Range("B18").Value = 400124
$ws.Range("C18").Value = 2344
$ws.Range("D18").Value = 332330
$ws.Range("E18").Value = 58004
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 55
$ws.Range("H18").Value = 9790

# Alemania
$ws.Range("B25").Value = 322194
$ws.Range("C25").Value = 1716
$ws.Range("D25").Value = 273500
$ws.Range("E25").Value = 39004
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9690

# Israel
$ws.Range("B27").Value = 289799
$ws.Range("C27").Value = 1941
$ws.Range("D27").Value = 225718
$ws.Range("E27").Value = 62167
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 28
$ws.Range("H27").Value = 1914

# Marruecos (now row 32)
$ws.Range("B32").Value = 149841
$ws.Range("C32").Value = 3443
$ws.Range("D32").Value = 124854
$ws.Range("E32").Value = 22415
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 42
$ws.Range("H32").Value = 2572

# Belgica (now row 33, unchanged numbers)
$ws.Range("B33").Value = 148981
$ws.Range("C33").Value = 5385
$ws.Range("D33").Value = 20072
$ws.Range("E33").Value = 118758
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = 10151

# Ecuador (now row 34, unchanged numbers)
$ws.Range("B34").Value = 146828
$ws.Range("C34").Value = 980
$ws.Range("D34").Value = 120511
$ws.Range("E34").Value = 14129
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 13
$ws.Range("H34").Value = 12188

# Libano
$ws.Range("B66").Value = 52558
$ws.Range("C66").Value = 1388
$ws.Range("D66").Value = 22719
$ws.Range("E66").Value = 29384
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 455

# Irlanda (now row 72)
$ws.Range("B72").Value = 41714
$ws.Range("C72").Value = 1011
$ws.Range("D72").Value = 23364
$ws.Range("E72").Value = 16526
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 1824

# Libia (now row 73, unchanged numbers)
$ws.Range("B73").Value = 41686
$ws.Range("C73").Value = 318
$ws.Range("D73").Value = 23791
$ws.Range("E73").Value = 17272
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 623

# Kenia (now row 74, unchanged numbers)
$ws.Range("B74").Value = 41158
$ws.Range("C74").Value = 538
$ws.Range("D74").Value = 31876
$ws.Range("E74").Value = 8522
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 760

# Republica de Macedonia
$ws.Range("B89").Value = 20555
$ws.Range("C89").Value = 392
$ws.Range("D89").Value = 16099
$ws.Range("E89").Value = 3671
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 785

# Curazao
$ws.Range("B174").Value = 571
$ws.Range("C174").Value = 22
$ws.Range("D174").Value = 308
$ws.Range("E174").Value = 262
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 1
